$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (shifts G..AC right to H..AD).
$ws.Columns("G").Insert()

# Give the new column G the same width as the old "Notes" column (now H),
# so Excel renders the two ex-adjacent 16-wide columns consistently.
$ws.Columns("G").ColumnWidth = $ws.Columns("H").ColumnWidth

# Header for the new "Plug Height" column.
$ws.Range("G3").Value = "Plug Height "

# Every data row (4-24) gets "Low" in the new Plug Height column.
for ($r = 4; $r -le 24; $r++) {
    $ws.Cells.Item($r, 7).Value = "Low"
}

# Leave the selection where the author left it after the edit.
$null = $ws.Range("G25").Select()
